$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Cells whose new Price text is purely numeric-looking (e.g. "213.01") are
# entered with a leading apostrophe so Excel stores them as text (matching
# the source data, which keeps trailing zeros / exact formatting) instead of
# silently converting them to floating-point numbers.

$ws.Range("D2").Value = "28.430.28"
$ws.Range("E2").Value = "  +3.56%  "
$ws.Range("D3").Value = "1.589.23"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "'213.01"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").Value = "'24.39"
$ws.Range("E8").Value = "  +8.00%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").Value = "'0.0886"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "1.593.39"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "28.447.51"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "'63.07"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "'229.64"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "0.0₃0706"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'7.47"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'151.66"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "'15.21"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "1.400.34"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  -9.20%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +7.97%  "
$ws.Range("D40").Value = "'0.541"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "'5.60"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").Value = "'0.982"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'63.17"
$ws.Range("E46").Value = "  -1.78%  "
$ws.Range("D47").Value = "1.726.04"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "'87.23"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("E51").Value = "  -0.94%  "
